# Apply crypto price/volume updates per commit "Updated cryptos list on Thu Sep 21 11:52:56 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.769.62"
$ws.Range("D3").Value = "'1.595.79"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'210.38"
$ws.Range("E5").Value = "  -2.68%  "
$ws.Range("D6").Value = "'0.506"
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").Value = "'19.58"
$ws.Range("E10").Value = "  -2.65%  "
$ws.Range("D11").Value = "'0.0834"
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").Value = "'1.817.11"
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("D13").Value = "'1.606.84"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").Value = "'4.06"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Value = "'0.529"
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("D16").Value = "'26.756.37"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").Value = "'63.50"
$ws.Range("E17").Value = "  -3.30%  "
$ws.Range("D18").Value = "'0.0₃0729"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").Value = "'209.37"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("D21").Value = "'6.71"
$ws.Range("E21").Value = "  -1.72%  "
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("E23").Value = "  -6.59%  "
$ws.Range("D24").Value = "'8.86"
$ws.Range("E24").Value = "  -2.77%  "
$ws.Range("D25").Value = "'146.53"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("E26").Value = "  +1.42%  "
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").Value = "'0.112"
$ws.Range("E28").Value = "  -4.50%  "
$ws.Range("D29").Value = "'15.31"
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("D30").Value = "'0.0501"
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("E31").Value = "  -2.49%  "
$ws.Range("E32").Value = "  -2.75%  "
$ws.Range("D33").Value = "'0.676"
$ws.Range("E33").Value = "  +24.54%  "
$ws.Range("E34").Value = "  -2.14%  "
$ws.Range("D35").Value = "'1.310.62"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  -3.11%  "
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").Value = "'0.820"
$ws.Range("E39").Value = "  -2.85%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  -2.09%  "
$ws.Range("D42").Value = "'2.17"
$ws.Range("E42").Value = "  -4.10%  "
$ws.Range("D43").Value = "'5.28"
$ws.Range("D44").Value = "'62.71"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").Value = "'1.730.92"
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("D46").Value = "'89.06"
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("D47").Value = "'1.61"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("D48").Value = "'0.808"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0509"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0975"
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.06%  "
